# Auto - Update data with bot!
# Updates the "title" (column D) values for a set of rows to match the
# blog/author "name" already present in column A of that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    6  = "슈퍼짱짱"
    16 = "Wonju Seo"
    19 = "꼬낄콘의 분석일지"
    20 = "ai-creator"
    23 = "Be the only one"
    28 = "로봇이 아닙니다 "
    32 = "데이터과학 삼학년"
    39 = "deadNstreet"
    42 = "IT_notepad"
    43 = "동신한의 조재성"
    44 = "Engineer-Ladder"
    45 = "dive-into-ds"
    46 = "BioinformaticsAndMe"
    47 = "shinminyoung"
    51 = "bskyvsion"
}

foreach ($row in $updates.Keys) {
    $ws.Range("D$row").Value = $updates[$row]
}
